# ---------------------------------------------------------------------------
# edit.ps1
#
# Reproduces the "Add files via upload" commit for static/Stimul_6.xlsx:
#   * Every cell in column B (rows 2-101) gets an explicit hyperlink pointing
#     at the matching "https://github.com/.../%D0%9AN.jpg?raw=true" GitHub
#     raw-image URL, while the cell keeps showing the short "KN.jpg" label
#     instead of the full URL text.
#   * Row 68, which previously duplicated row 18s picture number (K52),
#     is corrected to reference the previously-missing "K51.jpg" stimulus,
#     and rows 69-101 shift down by one picture number to match.
#   * The worksheet selection/scroll position is updated to B18:B101.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a hyperlink to every data row of column B (rows 2-101). The insertion
# order below matches the order the relationships were originally added in,
# so that the generated rId1..rId100 numbering lines up with the target file.
# The 5th argument (TextToDisplay) keeps the cell showing the short image
# file name ("K<number>.jpg") rather than the full URL.
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A1.jpg?raw=true", "", "", "К1.jpg")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A12.jpg?raw=true", "", "", "К12.jpg")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A11.jpg?raw=true", "", "", "К11.jpg")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A10.jpg?raw=true", "", "", "К10.jpg")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A19.jpg?raw=true", "", "", "К19.jpg")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A18.jpg?raw=true", "", "", "К18.jpg")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A29.jpg?raw=true", "", "", "К29.jpg")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A28.jpg?raw=true", "", "", "К28.jpg")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A27.jpg?raw=true", "", "", "К27.jpg")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A26.jpg?raw=true", "", "", "К26.jpg")
$ws.Hyperlinks.Add($ws.Range("B13"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A38.jpg?raw=true", "", "", "К38.jpg")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A37.jpg?raw=true", "", "", "К37.jpg")
$ws.Hyperlinks.Add($ws.Range("B16"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A45.jpg?raw=true", "", "", "К45.jpg")
$ws.Hyperlinks.Add($ws.Range("B15"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A44.jpg?raw=true", "", "", "К44.jpg")
$ws.Hyperlinks.Add($ws.Range("B14"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A43.jpg?raw=true", "", "", "К43.jpg")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A50.jpg?raw=true", "", "", "К50.jpg")
$ws.Hyperlinks.Add($ws.Range("B68"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A51.jpg?raw=true", "", "", "К51.jpg")
$ws.Hyperlinks.Add($ws.Range("B18"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A52.jpg?raw=true", "", "", "К52.jpg")
$ws.Hyperlinks.Add($ws.Range("B66"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A100.jpg?raw=true", "", "", "К100.jpg")
$ws.Hyperlinks.Add($ws.Range("B67"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A101.jpg?raw=true", "", "", "К101.jpg")
$ws.Hyperlinks.Add($ws.Range("B65"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A99.jpg?raw=true", "", "", "К99.jpg")
$ws.Hyperlinks.Add($ws.Range("B64"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A98.jpg?raw=true", "", "", "К98.jpg")
$ws.Hyperlinks.Add($ws.Range("B63"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A97.jpg?raw=true", "", "", "К97.jpg")
$ws.Hyperlinks.Add($ws.Range("B62"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A96.jpg?raw=true", "", "", "К96.jpg")
$ws.Hyperlinks.Add($ws.Range("B61"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A95.jpg?raw=true", "", "", "К95.jpg")
$ws.Hyperlinks.Add($ws.Range("B60"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A94.jpg?raw=true", "", "", "К94.jpg")
$ws.Hyperlinks.Add($ws.Range("B59"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A93.jpg?raw=true", "", "", "К93.jpg")
$ws.Hyperlinks.Add($ws.Range("B58"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A92.jpg?raw=true", "", "", "К92.jpg")
$ws.Hyperlinks.Add($ws.Range("B57"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A91.jpg?raw=true", "", "", "К91.jpg")
$ws.Hyperlinks.Add($ws.Range("B56"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A90.jpg?raw=true", "", "", "К90.jpg")
$ws.Hyperlinks.Add($ws.Range("B55"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A89.jpg?raw=true", "", "", "К89.jpg")
$ws.Hyperlinks.Add($ws.Range("B54"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A88.jpg?raw=true", "", "", "К88.jpg")
$ws.Hyperlinks.Add($ws.Range("B53"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A87.jpg?raw=true", "", "", "К87.jpg")
$ws.Hyperlinks.Add($ws.Range("B52"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A86.jpg?raw=true", "", "", "К86.jpg")
$ws.Hyperlinks.Add($ws.Range("B51"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A85.jpg?raw=true", "", "", "К85.jpg")
$ws.Hyperlinks.Add($ws.Range("B50"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A84.jpg?raw=true", "", "", "К84.jpg")
$ws.Hyperlinks.Add($ws.Range("B49"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A83.jpg?raw=true", "", "", "К83.jpg")
$ws.Hyperlinks.Add($ws.Range("B48"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A82.jpg?raw=true", "", "", "К82.jpg")
$ws.Hyperlinks.Add($ws.Range("B47"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A81.jpg?raw=true", "", "", "К81.jpg")
$ws.Hyperlinks.Add($ws.Range("B46"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A80.jpg?raw=true", "", "", "К80.jpg")
$ws.Hyperlinks.Add($ws.Range("B45"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A79.jpg?raw=true", "", "", "К79.jpg")
$ws.Hyperlinks.Add($ws.Range("B44"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A78.jpg?raw=true", "", "", "К78.jpg")
$ws.Hyperlinks.Add($ws.Range("B43"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A77.jpg?raw=true", "", "", "К77.jpg")
$ws.Hyperlinks.Add($ws.Range("B42"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A76.jpg?raw=true", "", "", "К76.jpg")
$ws.Hyperlinks.Add($ws.Range("B41"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A75.jpg?raw=true", "", "", "К75.jpg")
$ws.Hyperlinks.Add($ws.Range("B40"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A74.jpg?raw=true", "", "", "К74.jpg")
$ws.Hyperlinks.Add($ws.Range("B39"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A73.jpg?raw=true", "", "", "К73.jpg")
$ws.Hyperlinks.Add($ws.Range("B38"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A72.jpg?raw=true", "", "", "К72.jpg")
$ws.Hyperlinks.Add($ws.Range("B37"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A71.jpg?raw=true", "", "", "К71.jpg")
$ws.Hyperlinks.Add($ws.Range("B36"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A70.jpg?raw=true", "", "", "К70.jpg")
$ws.Hyperlinks.Add($ws.Range("B35"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A69.jpg?raw=true", "", "", "К69.jpg")
$ws.Hyperlinks.Add($ws.Range("B34"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A68.jpg?raw=true", "", "", "К68.jpg")
$ws.Hyperlinks.Add($ws.Range("B33"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A67.jpg?raw=true", "", "", "К67.jpg")
$ws.Hyperlinks.Add($ws.Range("B32"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A66.jpg?raw=true", "", "", "К66.jpg")
$ws.Hyperlinks.Add($ws.Range("B31"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A65.jpg?raw=true", "", "", "К65.jpg")
$ws.Hyperlinks.Add($ws.Range("B30"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A64.jpg?raw=true", "", "", "К64.jpg")
$ws.Hyperlinks.Add($ws.Range("B29"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A63.jpg?raw=true", "", "", "К63.jpg")
$ws.Hyperlinks.Add($ws.Range("B28"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A62.jpg?raw=true", "", "", "К62.jpg")
$ws.Hyperlinks.Add($ws.Range("B27"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A61.jpg?raw=true", "", "", "К61.jpg")
$ws.Hyperlinks.Add($ws.Range("B26"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A60.jpg?raw=true", "", "", "К60.jpg")
$ws.Hyperlinks.Add($ws.Range("B25"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A59.jpg?raw=true", "", "", "К59.jpg")
$ws.Hyperlinks.Add($ws.Range("B24"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A58.jpg?raw=true", "", "", "К58.jpg")
$ws.Hyperlinks.Add($ws.Range("B23"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A57.jpg?raw=true", "", "", "К57.jpg")
$ws.Hyperlinks.Add($ws.Range("B22"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A56.jpg?raw=true", "", "", "К56.jpg")
$ws.Hyperlinks.Add($ws.Range("B21"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A55.jpg?raw=true", "", "", "К55.jpg")
$ws.Hyperlinks.Add($ws.Range("B20"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A54.jpg?raw=true", "", "", "К54.jpg")
$ws.Hyperlinks.Add($ws.Range("B19"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A53.jpg?raw=true", "", "", "К53.jpg")
$ws.Hyperlinks.Add($ws.Range("B69"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A52.jpg?raw=true", "", "", "К52.jpg")
$ws.Hyperlinks.Add($ws.Range("B77"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A60.jpg?raw=true", "", "", "К60.jpg")
$ws.Hyperlinks.Add($ws.Range("B76"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A59.jpg?raw=true", "", "", "К59.jpg")
$ws.Hyperlinks.Add($ws.Range("B75"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A58.jpg?raw=true", "", "", "К58.jpg")
$ws.Hyperlinks.Add($ws.Range("B74"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A57.jpg?raw=true", "", "", "К57.jpg")
$ws.Hyperlinks.Add($ws.Range("B73"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A56.jpg?raw=true", "", "", "К56.jpg")
$ws.Hyperlinks.Add($ws.Range("B72"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A55.jpg?raw=true", "", "", "К55.jpg")
$ws.Hyperlinks.Add($ws.Range("B71"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A54.jpg?raw=true", "", "", "К54.jpg")
$ws.Hyperlinks.Add($ws.Range("B70"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A53.jpg?raw=true", "", "", "К53.jpg")
$ws.Hyperlinks.Add($ws.Range("B101"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A84.jpg?raw=true", "", "", "К84.jpg")
$ws.Hyperlinks.Add($ws.Range("B100"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A83.jpg?raw=true", "", "", "К83.jpg")
$ws.Hyperlinks.Add($ws.Range("B99"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A82.jpg?raw=true", "", "", "К82.jpg")
$ws.Hyperlinks.Add($ws.Range("B98"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A81.jpg?raw=true", "", "", "К81.jpg")
$ws.Hyperlinks.Add($ws.Range("B97"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A80.jpg?raw=true", "", "", "К80.jpg")
$ws.Hyperlinks.Add($ws.Range("B96"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A79.jpg?raw=true", "", "", "К79.jpg")
$ws.Hyperlinks.Add($ws.Range("B95"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A78.jpg?raw=true", "", "", "К78.jpg")
$ws.Hyperlinks.Add($ws.Range("B94"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A77.jpg?raw=true", "", "", "К77.jpg")
$ws.Hyperlinks.Add($ws.Range("B93"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A76.jpg?raw=true", "", "", "К76.jpg")
$ws.Hyperlinks.Add($ws.Range("B92"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A75.jpg?raw=true", "", "", "К75.jpg")
$ws.Hyperlinks.Add($ws.Range("B91"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A74.jpg?raw=true", "", "", "К74.jpg")
$ws.Hyperlinks.Add($ws.Range("B90"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A73.jpg?raw=true", "", "", "К73.jpg")
$ws.Hyperlinks.Add($ws.Range("B89"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A72.jpg?raw=true", "", "", "К72.jpg")
$ws.Hyperlinks.Add($ws.Range("B88"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A71.jpg?raw=true", "", "", "К71.jpg")
$ws.Hyperlinks.Add($ws.Range("B87"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A70.jpg?raw=true", "", "", "К70.jpg")
$ws.Hyperlinks.Add($ws.Range("B86"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A69.jpg?raw=true", "", "", "К69.jpg")
$ws.Hyperlinks.Add($ws.Range("B85"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A68.jpg?raw=true", "", "", "К68.jpg")
$ws.Hyperlinks.Add($ws.Range("B84"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A67.jpg?raw=true", "", "", "К67.jpg")
$ws.Hyperlinks.Add($ws.Range("B83"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A66.jpg?raw=true", "", "", "К66.jpg")
$ws.Hyperlinks.Add($ws.Range("B82"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A65.jpg?raw=true", "", "", "К65.jpg")
$ws.Hyperlinks.Add($ws.Range("B81"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A64.jpg?raw=true", "", "", "К64.jpg")
$ws.Hyperlinks.Add($ws.Range("B80"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A63.jpg?raw=true", "", "", "К63.jpg")
$ws.Hyperlinks.Add($ws.Range("B79"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A62.jpg?raw=true", "", "", "К62.jpg")
$ws.Hyperlinks.Add($ws.Range("B78"), "https://github.com/nastycify/ColorDotsExperiment5/blob/master/%D0%9A61.jpg?raw=true", "", "", "К61.jpg")

# Hyperlinks.Add() always stamps a brand-new cell style onto the target cell,
# even when the cell was already formatted with the built-in "Hyperlink" named
# style. Re-applying the same named style collapses the cell back onto the
# single shared "Hyperlink" style (as in the source workbook) instead of
# leaving behind 100 duplicate style records.
for ($r = 2; $r -le 101; $r++) {
    $cell = $ws.Range("B$r")
    $styleName = $cell.Style.Name
    $cell.Style = $styleName
}

# Update the sheet selection/scroll position to match the author's final view:
# the visible window was scrolled so row 86 is at the top, with B18:B101
# selected (active cell B18).
$win = $excel.ActiveWindow
$win.ScrollRow = 86
$win.ScrollColumn = 1
$ws.Range("B18:B101").Select()
